$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-11 Saturday" "2023-11-12 Sunday"

Replace-Text "32×26=" "36×49="
Replace-Text "22×55=" "21×30="
Replace-Text "51×35=" "27×60="
Replace-Text "15×52=" "73×99="
Replace-Text "96×31=" "67×77="
Replace-Text "99×86=" "50×14="
Replace-Text "13×38=" "20×99="
Replace-Text "67×48=" "28×79="
Replace-Text "50×85=" "62×85="
Replace-Text "63×21=" "86×48="
Replace-Text "25×64=" "49×81="
Replace-Text "83×33=" "45×61="
Replace-Text "50×51=" "29×18="
Replace-Text "67×20=" "58×95="
Replace-Text "11×52=" "69×22="
Replace-Text "45×28=" "55×36="
Replace-Text "26×29=" "53×81="
Replace-Text "39×51=" "68×36="
Replace-Text "14×47=" "54×42="
Replace-Text "92×60=" "78×70="
Replace-Text "79×78=" "49×41="
Replace-Text "81×71=" "12×65="
Replace-Text "31×16=" "14×36="
Replace-Text "30×68=" "64×65="
Replace-Text "87×69=" "37×63="
